# Swap the contents of column C (codeforiati:group-name) and column D
# (codeforiati:group-code) across the whole used range, including the
# header row. This mirrors the upstream codelist change where the
# "group-code" column now appears before the "group-name" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
